$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New property column "max_taxa" added after the existing "taxa_barstacks"
# column (BY). It gets its own header/description/value cells in column BZ,
# mirroring the formatting of the BY column (copy format first, then set
# the value so the numbers/styles/fonts/fills all line up).

# Row 2 ("name" row) -> BZ2 = "max_taxa"  (added to the shared-string table
# before the long description, so it lands on the lower string index).
$ws.Range("BY2").Copy($ws.Range("BZ2"))
$ws.Range("BZ2").Value2 = "max_taxa"

# Row 1 (long description row) -> BZ1
$ws.Range("BY1").Copy($ws.Range("BZ1"))
$ws.Range("BZ1").Value2 = "This determines how many taxaonmic names will be displayed in each graph at a maximum. Remaining entries are combined into 'others'."

# Data rows 3-7 -> BZ3:BZ7 = 18 (numeric, not a shared string)
foreach ($r in 3..7) {
  $src = $ws.Range("BY$r")
  $dst = $ws.Range("BZ$r")
  $src.Copy($dst)
  $dst.Value2 = 18
}

# Widen the new column so its legend/content isn't clipped.
$ws.Columns.Item(78).ColumnWidth = 17.3

# Restore the user's on-screen selection (was BU17, now shifted to the new
# last column BZ17).
$ws.Range("BZ17").Select() | Out-Null
